$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clean up header labels: drop the space so they read as single words
# (e.g. "First Name" -> "FirstName", "Last Name" -> "LastName")
$ws.Range("C1").Value = "FirstName"
$ws.Range("D1").Value = "LastName"

# Fix the mis-entered item name in row 4 ("iphone" -> "Flowers")
$ws.Range("A4").Value = "Flowers"

# Move the active selection to D1
$ws.Range("D1").Select()
